$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B11").Value = "空档区有空位时可以发动：将1张手牌放入空档区，然后抽2张牌。Test"
$ws.Range("C11").Style = $ws.Range("C2").Style

$ws.Range("B11").Select()
